# Add a "2022-Q4" quarterly sheet (right after "总计", before "2022-Q3")
# and a new summary row on "总计" for it.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by copying "2022-Q3" (same header
#    style/borders as all the other per-quarter sheets) to right
#    before itself, then renaming the copy and overwriting its data.
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($template)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# Drop the 4th (now stale) data row - 2022-Q4 only has 3 fund rows.
$newSheet.Rows.Item(5).Delete()

# Row 2
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'001628"
$newSheet.Cells.Item(2,3).Value = "招商体育文化休闲股票A"
$newSheet.Cells.Item(2,4).Value = "'2.33"
$newSheet.Cells.Item(2,5).Value = "'93.03"
$newSheet.Cells.Item(2,6).Value = "'5.09"
$newSheet.Cells.Item(2,7).Value = "'0.1186"
$newSheet.Cells.Item(2,8).Value = 5

# Row 3
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'003397"
$newSheet.Cells.Item(3,3).Value = "银华体育文化灵活配置混合"
$newSheet.Cells.Item(3,4).Value = "'0.53"
$newSheet.Cells.Item(3,5).Value = "'87.74"
$newSheet.Cells.Item(3,6).Value = "'3.56"
$newSheet.Cells.Item(3,7).Value = "'0.0189"
$newSheet.Cells.Item(3,8).Value = 7

# Row 4
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "'015395"
$newSheet.Cells.Item(4,3).Value = "招商体育文化休闲股票C"
$newSheet.Cells.Item(4,4).Value = "'0.29"
$newSheet.Cells.Item(4,5).Value = "'93.03"
$newSheet.Cells.Item(4,6).Value = "'5.09"
$newSheet.Cells.Item(4,7).Value = "'0.0148"
$newSheet.Cells.Item(4,8).Value = 5

# ------------------------------------------------------------------
# 2. Insert the new "2022-Q4" row into "总计" (right under the header)
#    and keep the running index (column A) in sync for every row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Re-apply the plain (unstyled) look of the data rows to the freshly
# inserted row - Insert() otherwise drags in a neighbouring style.
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)
$summary.Range("B3:D3").Copy()
$summary.Range("B2:D2").PasteSpecial(-4122)

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 3
$summary.Cells.Item(2,4).Value = 0.15

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(8,1).Value = 6
$summary.Cells.Item(9,1).Value = 7
$summary.Cells.Item(10,1).Value = 8

# ------------------------------------------------------------------
# 3. Restore the original tab-selection state ("2020-Q4" was the
#    active sheet before the edit; Copy() above moved focus away).
# ------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
